# Generate Report for Archive
#
# Updates the localization-status report:
#   - every cell that reads "Ready for handoff" now reads "In Translation"
#     (Overview!E2 + F2, zh-cn!C2, de-de!C2 — all share the same string)
#   - the now-narrower "In Translation" status column is re-sized to fit
#     (Overview columns E & F, and column C on each language sheet)

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth = 12.5

$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Range("E1:F1").Columns.ColumnWidth = $newWidth

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C1").ColumnWidth = $newWidth

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C1").ColumnWidth = $newWidth
